# "added images chart png"
# Repurpose the Blad3 (sheet3) mini-table + its bar chart from the old
# gzip/brotli-per-asset breakdown into a single-row "images before/after"
# comparison (original jpeg/png vs WebP), in MB instead of KB.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad3")

# --- Header row (B1/C1) ---------------------------------------------------
$ws.Range("B1").Value = "Original (jpeg and png)"
$ws.Range("C1").Value = "WebP version"

# --- Data row 2: the single remaining data point --------------------------
$ws.Range("A2").Value = "Total size of 7 images"
$ws.Range("B2").Value = 4.43
$ws.Range("C2").Value = 1.27

# New custom number format (MB, one-or-two decimals) for the two new values.
$ws.Range("B2:C2").NumberFormat = '0.0#\ "MB"'

# --- Rows 3 & 4 no longer hold data: drop the old per-asset rows ----------
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("B3:C4").ClearContents()

# --- Page setup now configured for this sheet (portrait / A4-ish letter) --
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Chart: now sourced from Blad3 itself instead of Blad1 ----------------
$chart = $ws.ChartObjects(1).Chart

$s1 = $chart.SeriesCollection(1)
$s1.Name = "=Blad3!$B$1"
$s1.Values = "=Blad3!$B$2"
$s1.XValues = "=Blad3!$A$2"

$s2 = $chart.SeriesCollection(2)
$s2.Name = "=Blad3!$C$1"
$s2.Values = "=Blad3!$C$2"
$s2.XValues = "=Blad3!$A$2"

$chart.Axes(2).TickLabels.NumberFormat = '0.0# "MB"'

# --- Selection cursor, as saved in the authored workbook (must be last) ---
$ws.Range("C14").Select()
